$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains exact text formatting (e.g. trailing zeros,
# thousand-separator dots) instead of being auto-converted to numbers.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '42.360.14'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '2.281.28'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.49%  '
$ws.Range('D5').Value = '310.11'
$ws.Range('E5').Value = '  -3.99%  '
$ws.Range('D6').Value = '103.31'
$ws.Range('E6').Value = '  +0.95%  '
$ws.Range('D7').Value = '0.620'
$ws.Range('E7').Value = '  -1.35%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.23%  '
$ws.Range('D9').Value = '0.599'
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('D10').Value = '38.56'
$ws.Range('E10').Value = '  -2.70%  '
$ws.Range('D11').Value = '0.0893'
$ws.Range('E11').Value = '  -1.33%  '
$ws.Range('D12').Value = '8.20'
$ws.Range('E12').Value = '  -1.21%  '
$ws.Range('E13').Value = '  +0.85%  '
$ws.Range('D14').Value = '0.970'
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('D15').Value = '15.00'
$ws.Range('E15').Value = '  -0.41%  '
$ws.Range('D16').Value = '2.627.10'
$ws.Range('E16').Value = '  -0.22%  '
$ws.Range('D17').Value = '2.282.34'
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = '42.326.86'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').Value = '7.25'
$ws.Range('E19').Value = '  -1.24%  '
$ws.Range('E20').Value = '  -1.32%  '
$ws.Range('D21').Value = '12.97'
$ws.Range('E21').Value = '  +2.09%  '
$ws.Range('D22').Value = '72.64'
$ws.Range('E22').Value = '  -0.50%  '
$ws.Range('D23').Value = '3.39'
$ws.Range('E23').Value = '  -7.26%  '
$ws.Range('D24').Value = '262.09'
$ws.Range('E24').Value = '  -2.11%  '
$ws.Range('D25').Value = '2.17'
$ws.Range('E25').Value = '  -2.52%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('D27').Value = '10.62'
$ws.Range('E27').Value = '  -1.95%  '
$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D28').Value = '6.87'
$ws.Range('E28').Value = '  +13.51%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.28'
$ws.Range('E29').Value = '  -1.54%  '
$ws.Range('D30').Value = '22.08'
$ws.Range('E30').Value = '  -1.56%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = '35.62'
$ws.Range('E31').Value = '  -6.30%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').Value = '164.46'
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('D33').Value = '0.0851'
$ws.Range('E33').Value = '  -2.77%  '
$ws.Range('E34').Value = '  -2.93%  '
$ws.Range('E35').Value = '  +0.36%  '
$ws.Range('E36').Value = '  -4.24%  '
$ws.Range('D37').Value = '4.48'
$ws.Range('E37').Value = '  -2.47%  '
$ws.Range('D38').Value = '0.0346'
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('D39').Value = '3.67'
$ws.Range('E39').Value = '  -0.29%  '
$ws.Range('D40').Value = '2.69'
$ws.Range('E40').Value = '  -1.98%  '
$ws.Range('E41').Value = '  +1.95%  '
$ws.Range('D42').Value = '98.20'
$ws.Range('E42').Value = '  +8.32%  '
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').Value = '68.38'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('D45').Value = '0.225'
$ws.Range('E45').Value = '  -0.33%  '
$ws.Range('D46').Value = '11.88'
$ws.Range('E46').Value = '  -2.77%  '
$ws.Range('D47').Value = '1.701.83'
$ws.Range('E47').Value = '  +6.97%  '
$ws.Range('D48').Value = '109.82'
$ws.Range('E48').Value = '  -2.66%  '
$ws.Range('D49').Value = '76.62'
$ws.Range('E49').Value = '  -4.11%  '
$ws.Range('D50').Value = '8.59'
$ws.Range('E50').Value = '  -3.74%  '
$ws.Range('E51').Value = '  -2.46%  '
